$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings are not
# reinterpreted as numbers by Excel when the values are assigned below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.266.04"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.550.26"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "208.98"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "0.484"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "23.34"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "0.0583"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.771.56"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "1.548.89"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "28.263.93"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").Value = "60.37"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("D18").Value = "225.91"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "7.29"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "0.0₃0673"
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "3.91"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D25").Value = "147.69"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").Value = "14.75"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").Value = "0.0464"
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "1.380.53"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").Value = "2.58"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "0.509"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.773"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.0464"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "5.40"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "61.55"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D48").Value = "1.685.29"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "85.27"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "41.29"
$ws.Range("E51").Value = "  +3.64%  "

# Restore the original (default) cell style now that the text values are set,
# so no new/explicit number-format style lingers on these cells.
$ws.Range("D2:D51").Style = "Normal"
